$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add a new worksheet "Sheet2" positioned right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 with header + old data row + new data row
$ws2.Range("A1").Value = "CREDIT.ACCT.NO"
$ws2.Range("B1").Value = "CREDIT.THEIR.REF"
$ws2.Range("A2").Value = 1000140984
$ws2.Range("B2").Value = "LDA0610297"
$ws2.Range("A3").Value = 1005667626
$ws2.Range("B3").Value = "CLK0601335"

# Select whole sheet on Sheet2 (no active cell override, default selection)
$ws2.Cells.Select()

# Update Sheet1 row 2 with new data
$ws1.Range("A2").Value = 1005667626
$ws1.Range("B2").Value = "CLK0601335"

# Set selection on Sheet1 to A2 and make Sheet1 active/tab-selected
$ws1.Activate()
$ws1.Range("A2").Select()
